$d = $word.ActiveDocument

$p = $d.Paragraphs(2)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Deleting all fields in</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> the table </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>accessory_received_quantity</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>:</w:t></w:r></w:p>')

$p = $d.Paragraphs(3)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:t>delete</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> from </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>accessory_received_quantity</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

$p = $d.Paragraphs(4)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:t>where</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>received_quantity_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> &lt; 1000</w:t></w:r></w:p>')

$p = $d.Paragraphs(5)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Resetting all fields in the table: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>accessory_inventory</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>:</w:t></w:r></w:p>')

$p = $d.Paragraphs(6)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>update</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>accessory_inventory</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

$p = $d.Paragraphs(7)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="720"/></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:t>set</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>current_inventory</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = 0, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sales_quantity</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = 0, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sales_amount</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = 0, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>purchased_amount</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = 0, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>total_current_inventory</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = 0, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>total_sales_quantity</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = 0, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>total_sales_amount</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = 0, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>total_purchased_amount</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = 0, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>current_inventory_amount</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = 0, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>total_purchased_quantity</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = 0, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>purchased_quantity</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = 0, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>received_date</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = null</w:t></w:r></w:p>')

$p = $d.Paragraphs(8)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="720"/></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:t>where</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>accessory_inventory_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> &lt; 1000</w:t></w:r></w:p>')

$p = $d.Paragraphs(10)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Deleting all fields in the table </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>phones</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>_received_quantity</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>:</w:t></w:r></w:p>')

$p = $d.Paragraphs(11)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:t>delete</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> from </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>phones_received_quantity</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

$p = $d.Paragraphs(12)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:t>where</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>received_quantity_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> &lt; 1000</w:t></w:r></w:p>')

$p = $d.Paragraphs(13)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Resetting all fields in the table: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>unlocked_phones_inventory</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>:</w:t></w:r></w:p>')

$p = $d.Paragraphs(14)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:t>update</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>unlocked_phones_inventory</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

$p = $d.Paragraphs(15)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:t xml:space="preserve">set </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>current_inventory</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = 0, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sales_quantity</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = 0, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sales_amount</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = 0, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>purchased_amount</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = 0, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>total_current_inventory</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = 0, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>total_sales_quantity</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = 0, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>total_sales_amount</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = 0, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>total_purchased_amount</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = 0, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>current_inventory_amount</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = 0, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>total_purchased_quantity</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = 0, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>purchased_quantity</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = 0, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>received_date</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = null</w:t></w:r></w:p>')

$p = $d.Paragraphs(16)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="720"/></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:t>where</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>unlocked_phone</w:t></w:r><w:r><w:t>_inventory_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> &lt; 1000</w:t></w:r></w:p>')
